# Insert a new row at position 104 (shifting existing rows 104:199 down to 105:200)
# and populate it with the new data record, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("104:104").Insert()

$ws.Cells.Item(104,1).Value2  = 1
$ws.Cells.Item(104,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(104,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(104,4).Value2  = 44484
$ws.Cells.Item(104,5).Value2  = 15
$ws.Cells.Item(104,6).Value2  = 100114013
$ws.Cells.Item(104,7).Value2  = "Zanahoria"
$ws.Cells.Item(104,8).Value2  = "Sin especificar"
$ws.Cells.Item(104,9).Value2  = "Primera"
$ws.Cells.Item(104,10).Value2 = 70
$ws.Cells.Item(104,11).Value2 = 11000
$ws.Cells.Item(104,12).Value2 = 12000
$ws.Cells.Item(104,13).Value2 = 11500
$ws.Cells.Item(104,14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(104,15).Value2 = "Valle de Camiña"
$ws.Cells.Item(104,16).Value2 = 460
$ws.Cells.Item(104,17).Value2 = 25
$ws.Cells.Item(104,18).Value2 = "Hortaliza"
